$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently has 8 columns: Course, Teacher, Discipline, Venue,
# Status, LateIn, Left Early, Remarks (A1:H19). We need to insert a new
# "Session" column before "LateIn" (which itself gets renamed "Late in"),
# growing the table to 9 columns (A1:I19):
#   Course, Teacher, Discipline, Venue, Status, Session, Late in, Left Early, Remarks
$lo = $ws.ListObjects.Item(1)

# Shift the existing header labels right by one column (F->G, G->H, H->I)
# before growing the table, so the new slot opens up at F.
$ws.Range("H1").Value = $ws.Range("G1").Value()
$ws.Range("G1").Value = $ws.Range("F1").Value()
$ws.Range("F1").Value = "Session"

# Grow the table range to include the new column I (the engine will assign
# a placeholder name for the brand-new column; we fix that up below once
# the real header text is in place).
$lo.Resize($ws.Range("A1:I19"))

# Put "Remarks" into the newly-added column I now that the table already
# spans it, and rename the relocated "LateIn" header to "Late in".
$ws.Range("I1").Value = "Remarks"
$ws.Range("G1").Value = "Late in"

# Slightly widen column E to match the post-edit layout.
$ws.Columns.Item(5).ColumnWidth = 14.5
